$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$quot = [char]8220   # “
$rquot = [char]8221  # ”
$ndash = [char]8211  # –

# ---------------------------------------------------------------------------
# 1) "Adding ---- ... " paragraph: split the single run into three runs so
#    the quoted "git" is wrapped with gramStart/gramEnd proofErr markers.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Adding ---- " + $quot + " git add filename1 filename2" + $rquot + " ---- separate files with spaces to add multiple at once. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $addingXml = $pkgOpen + `
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
        '<w:r><w:t xml:space="preserve">Adding ---- </w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:t>&#8220; git</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:t xml:space="preserve"> add filename1 filename2&#8221; ---- separate files with spaces to add multiple at once. </w:t></w:r>' + `
        '</w:p>' + $pkgClose
    $rng.Delete()
    $rng.InsertXML($addingXml)
}

# ---------------------------------------------------------------------------
# 2) The empty sub-bullet paragraph right after "Adding ----" gets two runs
#    describing the "git commit" command.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Adding ---- " + $quot + " git add filename1 filename2" + $rquot + " ---- separate files with spaces to add multiple at once. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found2) {
    $rng2.Collapse(0)
    $emptyPara2 = $rng2.Next(4, 1)
    $commitXml = $pkgOpen + `
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
        '<w:r><w:t xml:space="preserve">Committing &#8211; use &#8220;git commit </w:t></w:r>' + `
        '<w:r><w:t>-m &#8220;Message here in quotes&#8221; &#8220; &#8211; always  use git commit -m &#8220;My message&#8221; ----- while you are committing new changes.</w:t></w:r>' + `
        '</w:p>' + $pkgClose
    $emptyPara2.InsertXML($commitXml)
}

# ---------------------------------------------------------------------------
# 3) The trailing empty paragraph (after "Repository ---") gets a run, and
#    two brand-new sub-bullet paragraphs are appended right after it.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Repository --- actual git repository " + "the .git" + " folder where you are making the changes or commits", `
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found3) {
    $rng3.Collapse(0)
    $emptyPara3 = $rng3.Next(4, 1)
    $emptyPara3.InsertAfter("Git Commit -m " + $quot + "my message'")

    $afterGitCommit = $d.Range($emptyPara3.End, $emptyPara3.End)
    $newParasXml = $pkgOpen + `
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
        '<w:r><w:t xml:space="preserve">We use the &#8220;git commit&#8221; command to </w:t></w:r>' + `
        '<w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:t>actually commit</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:t xml:space="preserve"> changes from the staging area</w:t></w:r>' + `
        '</w:p>' + `
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
        '<w:r><w:t>When making a commit, we need to provide a commit message that summarized the changes and work snapshotted in the commit</w:t></w:r>' + `
        '</w:p>' + $pkgClose
    $afterGitCommit.InsertXML($newParasXml)
}
